$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-5 hold text-formatted numbers/percentages (e.g. "41.00", "-109.00",
# "27.33%") stored as plain text, not real numeric values. A leading
# apostrophe forces Excel to keep the literal text instead of re-parsing it
# as a number/percentage when assigned through .Value.

# New row 3 = Ochieng Charles's data (previously row 5)
$ws.Range("A3").Value = "Ochieng Charles"
$ws.Range("B3").Value = "'41.00"
$ws.Range("D3").Value = "'-109.00"
$ws.Range("E3").Value = "'27.33%"

# New row 4 = Lenah Cheloti's data (previously row 3)
$ws.Range("A4").Value = "Lenah Cheloti"
$ws.Range("B4").Value = "'30.00"
$ws.Range("D4").Value = "'-120.00"
$ws.Range("E4").Value = "'20.00%"

# New row 5 = Moses  Ngugi's data (previously row 4)
$ws.Range("A5").Value = "Moses  Ngugi"
$ws.Range("B5").Value = "'24.00"
$ws.Range("D5").Value = "'-126.00"
$ws.Range("E5").Value = "'16.00%"
